# Limit max searches displayed to 10
#
# - Drop the stray "_GoBack" bookmark that sat at the end of the
#   "Case-insensitive search" bullet.
# - Add a trailing space run to the "Usability" bullet.
# - Add two new sub-bullets ("Invalid commands caught" and
#   "Max 10 results displayed") under it, re-homing the "_GoBack"
#   bookmark onto the very end of the last new bullet.
#
# NOTE: paragraph objects are re-fetched from $d.Paragraphs by index
# after every mutation instead of being cached/chained via .Next,
# since stale paragraph handles do not reliably track the document
# after edits in this host.

$d = $word.ActiveDocument

# --- 1. Remove the old _GoBack bookmark -------------------------------
$oldGoBack = $d.Bookmarks.Item("_GoBack")
$oldGoBack.Delete()

# --- 2. Find the "Usability ... POODLE command line tool" paragraph ---
$usabilityIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*POODLE command line tool*") {
        $usabilityIdx = $i
        break
    }
}

# Append a trailing space run to that paragraph.
$usabilityPara = $d.Paragraphs.Item($usabilityIdx)
$usabilityPara.Range.InsertAfter(" ")

# --- 3. Insert "Invalid commands caught" as a new sub-bullet ----------
$usabilityPara = $d.Paragraphs.Item($usabilityIdx)
$usabilityPara.Range.InsertParagraphAfter()

$invalidIdx = $usabilityIdx + 1
$invalidPara = $d.Paragraphs.Item($invalidIdx)
$invalidPara.Range.ListFormat.ListLevelNumber = 2
$invalidPara.Range.Text = "Invalid commands caught"

# --- 4. Insert "Max 10 results displayed" as the next sub-bullet ------
$invalidPara = $d.Paragraphs.Item($invalidIdx)
$invalidPara.Range.InsertParagraphAfter()

$maxIdx = $invalidIdx + 1
$maxPara = $d.Paragraphs.Item($maxIdx)
$maxPara.Range.ListFormat.ListLevelNumber = 2

# Type the text with two throwaway trailing characters so that the
# bookmark can be planted at a safe, unambiguous mid-run offset (the
# true end-of-paragraph-content boundary is ambiguous and otherwise
# snaps to cover the whole paragraph). Once the bookmark is anchored,
# delete the placeholder characters after it so it ends up collapsed
# right at the end of the real text, matching the original "_GoBack"
# placement style.
$maxPara.Range.Text = "Max 10 results displayedXX"

$maxPara = $d.Paragraphs.Item($maxIdx)
$contentEnd = $maxPara.Range.End - 1
$bmPos = $contentEnd - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$newGoBack = $d.Bookmarks.Item("_GoBack")
$placeholderRange = $d.Range($newGoBack.End, $contentEnd)
$placeholderRange.Delete()

Write-Output "Done"
